$d = $word.ActiveDocument

# --- Edit 1 -----------------------------------------------------------
# "Petra (neutral raised_eyebrow): That's what you were talking about...?"
#   -> "Petra (neutral expressionless): That's what you were talking about...?"
$d.Content.Find.Execute(
    "Petra (neutral raised_eyebrow): That’s what you were talking about…?",
    $true, $false, $false, $false, $false, $true, 1, $false,
    "Petra (neutral expressionless): That’s what you were talking about…?",
    2) | Out-Null

# --- Edit 2a ------------------------------------------------------------
# "Petra (neutral raised_eyebrow): You don't play any sports..."
#   -> "Petra (neutral thinking): You don't play any sports..."
$d.Content.Find.Execute(
    "Petra (neutral raised_eyebrow): You don’t play any sports or anything, and Lilith’s pretty athletic…",
    $true, $false, $false, $false, $false, $true, 1, $false,
    "Petra (neutral thinking): You don’t play any sports or anything, and Lilith’s pretty athletic…",
    2) | Out-Null

# --- Edit 2b ------------------------------------------------------------
# "Petra (neutral thinking): But at the end of the day you "
#   -> "Petra (neutral expressionless): But at the end of the day you "
$d.Content.Find.Execute(
    "Petra (neutral thinking): But at the end of the day you ",
    $true, $false, $false, $false, $false, $true, 1, $false,
    "Petra (neutral expressionless): But at the end of the day you ",
    2) | Out-Null

# --- Edit 3 ---------------------------------------------------------------
# Merge the three runs "Pro: I'm probably gonna go home and finish my " +
# "career form" + "." into a single run of unbroken text. The combined
# text already reads correctly, so a direct same-text assignment is a
# no-op; write a placeholder first and then the real text so the runs
# actually get collapsed into one.
foreach ($p in $d.Paragraphs) {
    $t = $p.Range.Text
    if ($t -like "Pro: I*m probably gonna go home and finish my career form.*") {
        $start = $p.Range.Start
        $end = $p.Range.End - 1
        $r = $d.Range($start, $end)
        $r.Text = "@@TEMP@@"
        $r2 = $d.Range($start, $start + 8)
        $r2.Text = "Pro: I’m probably gonna go home and finish my career form."
        break
    }
}

# --- Edit 4 -----------------------------------------------------------
# Trailing space after the closing brace on the very last line of the
# document gets trimmed: "} " -> "}"
$d.Content.Find.Execute(
    "} ",
    $true, $false, $false, $false, $false, $true, 1, $false,
    "}",
    2) | Out-Null
